$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: F28/J28/K28 were stored as text; convert them to real numbers ---
$ws.Range("F28").Value = 9935558059
$ws.Range("J28").Value = 26.9107501
$ws.Range("K28").Value = 80.9465604

# --- New row 29: another "Pranjal Shukla" entry for the new listing ---
$ws.Range("E29").Value = "Pranjal Shukla"
# F29/J29/K29 keep the legacy text-typed formatting (leading apostrophe forces text,
# matching how this sheet's newest rows get appended before being normalized to numbers)
$ws.Range("F29").Value = "'9935558059"
$ws.Range("G29").Value = "propques"
$ws.Range("H29").Value = "pranjalshukla800@gmail.com"
$ws.Range("I29").Value = "67c15dea7a5466fe052ba655"
$ws.Range("J29").Value = "'26.9176756"
$ws.Range("K29").Value = "'80.9539367"
$ws.Range("L29").Value = "LU New Campus Road, Aliganj, Lucknow, Uttar Pradesh, 226026, India"
